$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "AlphaFiberF"

# Add new row 16, mirroring the pattern of row 15 (A=index, B=shared label, C:M = 1)
$row = 16

# Copy row 15's formatting (e.g. the bold/bordered style on column A) down to row 16
$ws.Range("A15").Copy($ws.Range("A16"))

$ws.Cells.Item($row, 1).Value = 14
$ws.Cells.Item($row, 2).Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 13; $col++) {
    $ws.Cells.Item($row, $col).Value = 1
}
